{"js": "// Split the single-run bibliography paragraph (\"FURTADO, N.; KAWAMOTO, E. ...\")\n// into the same run, but with a pair of manual line breaks (<w:br/><w:br/>)\n// inserted between each of the six reference entries that were previously\n// concatenated back-to-back with no separator.\n\nconst segments = [\n  \"FURTADO, N.; KAWAMOTO, E. Avalia\u00e7\u00e3o de Projetos de Transporte. S\u00e3o Carlos: Servi\u00e7o Gr\u00e1fico EESC-USP, 2002. 254 p.\",\n  \"POWER, D. J. Decision Support Systems. London: Quorum Books, 2002. 251 p.\",\n  \"GOMES, L. F. A. M.; GOMES, C. F. S.; ALMEIDA, A. T, Tomada de Decis\u00e3o Gerencial: enfoque multicrit\u00e9rio, S\u00e3o Paulo: Atlas, 2002.\",\n  \"SHIMIZU, T., Decis\u00e3o nas Organiza\u00e7\u00f5es: introdu\u00e7\u00e3o aos problemas de decis\u00e3o encontrados nas organiza\u00e7\u00f5es e nos sistemas de apoio \u00e0 decis\u00e3o, S\u00e3o Paulo: Atlas, 2001.\",\n  \"DEVLIN, G. (ed.). Decision Support Systems: advances in. Zagreb: Intech, 2010. 342 p.\",\n  \"GARC\u00cdA-D\u00cdAZ, V. Algorithms in Decision Support Systems. Basel: MDPI, 2020. 147 p.\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Find the paragraph that starts the bibliography entry (contains the first\n// segment) so the script doesn't depend on a hard-coded paragraph index.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  para.load(\"text\");\n}\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text && text.indexOf(segments[0]) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Bibliography paragraph not found\");\n}\n\n// Insert the two manual line breaks right after each segment (except the\n// last) by searching for the segment text inside the paragraph and inserting\n// \"after\" the matched range. Working front-to-back is safe because each\n// search is scoped to the (still single-run) paragraph and matches plain\n// text only once.\nfor (let i = 0; i < segments.length - 1; i++) {\n  const results = target.search(segments[i], { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Segment not found: \" + segments[i]);\n  }\n\n  results.items[0].insertText(\"\\u000b\\u000b\", \"After\");\n  await context.sync();\n}\n", "ps1": "# Split the single-run bibliography paragraph (\"FURTADO, N.; KAWAMOTO, E. ...\")\n# so that a pair of manual line breaks (<w:br/><w:br/>) separates each of the\n# six reference entries that were previously concatenated back-to-back with\n# no separator.\n\n$d = $word.ActiveDocument\n\n# Manual line break character (vbVerticalTab / Chr(11)) -- cast to [string]\n# first so \"+\" concatenates characters instead of adding their numeric values.\n$brk = [string][char]11\n\n$segments = @(\n  \"FURTADO, N.; KAWAMOTO, E. Avalia\u00e7\u00e3o de Projetos de Transporte. S\u00e3o Carlos: Servi\u00e7o Gr\u00e1fico EESC-USP, 2002. 254 p.\",\n  \"POWER, D. J. Decision Support Systems. London: Quorum Books, 2002. 251 p.\",\n  \"GOMES, L. F. A. M.; GOMES, C. F. S.; ALMEIDA, A. T, Tomada de Decis\u00e3o Gerencial: enfoque multicrit\u00e9rio, S\u00e3o Paulo: Atlas, 2002.\",\n  \"SHIMIZU, T., Decis\u00e3o nas Organiza\u00e7\u00f5es: introdu\u00e7\u00e3o aos problemas de decis\u00e3o encontrados nas organiza\u00e7\u00f5es e nos sistemas de apoio \u00e0 decis\u00e3o, S\u00e3o Paulo: Atlas, 2001.\",\n  \"DEVLIN, G. (ed.). Decision Support Systems: advances in. Zagreb: Intech, 2010. 342 p.\",\n  \"GARC\u00cdA-D\u00cdAZ, V. Algorithms in Decision Support Systems. Basel: MDPI, 2020. 147 p.\"\n)\n\n# Locate the bibliography paragraph by its distinctive leading text instead of\n# a hard-coded paragraph index.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith($segments[0])) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Bibliography paragraph not found\"\n}\n\n# Insert the two manual line breaks right after each segment (except the\n# last). Re-duplicating the paragraph range and re-running Find each time\n# keeps the search scoped correctly even as the paragraph grows.\nfor ($i = 0; $i -lt $segments.Length - 1; $i++) {\n    $r = $target.Range.Duplicate\n    $found = $r.Find.Execute($segments[$i])\n    if (-not $found) {\n        throw \"Segment not found: $($segments[$i])\"\n    }\n    $r.Collapse(0)  # wdCollapseEnd\n    $r.InsertAfter($brk + $brk)\n}\n"}
